# LOM3077.xlsx update
# - Delete old row 13 (the stray "471420 - Carlos Antonio Reis Pereira Baptista"
#   row with no column-A label), which shifts every row below it up by one.
# - Re-point a handful of B/C cells whose "paired" text no longer lines up
#   correctly after the shift, reusing existing text where possible and
#   adding the one genuinely new value ("Semestral").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture values we still need before they move / disappear.
$carlosName      = $ws.Range("B13").Value2   # "471420 - Carlos Antonio Reis Pereira Baptista"
$aulaExpositiva  = $ws.Range("B19").Value2   # "Aula expositiva ..." (Método: row, before shift)
$notaFinal       = $ws.Range("B20").Value2   # "A nota final sera..." (Critério: row, before shift)
$recuperacao     = $ws.Range("B21").Value2   # "A recuperação sera..." (Norma de recuperação row, before shift)

# Remove the stray row; everything below slides up (row 24 -> 23, etc.)
$ws.Rows("13").Delete()

# Objetivos: row now holds the professor info instead of the long text.
$ws.Range("B10").Value = $carlosName
$ws.Range("C10").Value = $carlosName

# Programa resumido: row (was row 14, now row 13) gets "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Programa: row (was row 16, now row 15) gets the activation date. Copy
# (rather than re-typing the string) so the date-looking text stays a
# shared string instead of being reparsed into a date serial, then fix up
# the pasted-in style (Copy also brings B8's style) to match each column.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("B8").Copy($ws.Range("C15"))
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# Método: row (was row 19, now row 18) gets the professor info again.
$ws.Range("B18").Value = $carlosName
$ws.Range("C18").Value = $carlosName

# Critério: row (was row 20, now row 19) gets the "Aula expositiva..." text.
$ws.Range("B19").Value = $aulaExpositiva
$ws.Range("C19").Value = $aulaExpositiva

# Norma de recuperação: row (was row 21, now row 20) gets the "nota final" text.
$ws.Range("B20").Value = $notaFinal
$ws.Range("C20").Value = $notaFinal

# Bibliografia: row (was row 22, now row 21) gets the "recuperação" text.
$ws.Range("B21").Value = $recuperacao
$ws.Range("C21").Value = $recuperacao
